$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.822.64"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").Value = "2.211.21"
$ws.Range("E3").Value = "  -1.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.24"
$ws.Range("E5").Value = "  +4.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.16"
$ws.Range("E7").Value = "  +2.51%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  -1.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.36"
$ws.Range("E10").Value = "  +2.55%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").Value = "  -2.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.03"
$ws.Range("E12").Value = "  +1.21%  "

# Row 13
$ws.Range("E13").Value = "  +1.19%  "

# Row 14
$ws.Range("D14").Value = "2.546.05"
$ws.Range("E14").Value = "  -1.54%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.45"
$ws.Range("E15").Value = "  -0.91%  "

# Row 16
$ws.Range("D16").Value = "2.208.11"
$ws.Range("E16").Value = "  -1.52%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.782"
$ws.Range("E17").Value = "  -1.41%  "

# Row 18
$ws.Range("D18").Value = "42.791.28"
$ws.Range("E18").Value = "  -0.54%  "

# Row 19
$ws.Range("E19").Value = "  -2.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.07"
$ws.Range("E20").Value = "  -0.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.96"
$ws.Range("E21").Value = "  -0.30%  "

# Row 22
$ws.Range("E22").Value = "  +2.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.19"
$ws.Range("E23").Value = "  +0.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.31"
$ws.Range("E24").Value = "  -6.05%  "

# Row 25
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "42.65"
$ws.Range("E26").Value = "  +11.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("E27").Value = "  -0.89%  "

# Row 28
$ws.Range("E28").Value = "  -2.32%  "

# Row 29
$ws.Range("E29").Value = "  -1.78%  "

# Row 30
$ws.Range("E30").Value = "  -0.34%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.93"
$ws.Range("E31").Value = "  +0.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.36"
$ws.Range("E32").Value = "  +0.31%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0874"
$ws.Range("E33").Value = "  +9.52%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.22"
$ws.Range("E34").Value = "  -1.51%  "

# Row 35
$ws.Range("E35").Value = "  -0.52%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0358"
$ws.Range("E36").Value = "  +8.16%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.107"
$ws.Range("E37").Value = "  -2.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.33"
$ws.Range("E38").Value = "  -0.70%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.02"
$ws.Range("E39").Value = "  -0.91%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("E40").Value = "  +16.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.11"
$ws.Range("E41").Value = "  -0.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.201"
$ws.Range("E42").Value = "  -2.40%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.31"
$ws.Range("E43").Value = "  -2.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.16"
$ws.Range("E44").Value = "  +1.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.08"
$ws.Range("E45").Value = "  -1.57%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.39"
$ws.Range("E46").Value = "  -3.67%  "

# Row 47
$ws.Range("B47").Value = "WOONetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.471"
$ws.Range("E47").Value = "  -2.85%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0976"
$ws.Range("E48").Value = "  -1.32%  "

# Row 49
$ws.Range("E49").Value = "  +0.64%  "

# Row 50
$ws.Range("E50").Value = "  -1.20%  "

# Row 51
$ws.Range("D51").Value = "2.431.62"
$ws.Range("E51").Value = "  -1.10%  "
